# Product Backlog.xlsx -- update user-story wording and fill in Story Points
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Re-word the User Story column (B2:B21), rows 2-21, row 12 unchanged ---
$ws.Range("B2").Value = "As a player, I want to spectate the game after death so that I can see the game’s outcome"
$ws.Range("B3").Value = "As a player, I want to see past effects on my character so that I can understand its current state"
$ws.Range("B4").Value = "As a player, I want to see recent actions so that I can understand others’ behavior"
$ws.Range("B5").Value = "As a player, I want to see a room’s effects so that I can decide if I want to enter"
$ws.Range("B6").Value = "As a player, I want to play a tutorial so that I can learn how to play the game."
$ws.Range("B7").Value = "As a player, I want to see my local gameplay stats (games played, times died, etc.) so that I can see, in number form, how my games tend to go"
$ws.Range("B8").Value = "As a player, I want to create a custom character (with stats) so that I can create a bit of investment into the world and game"
$ws.Range("B9").Value = "As a player, I want a place to view the character lore, so that I can immerse myself more in the world and lore "
$ws.Range("B10").Value = "As a player, I want the game elements to evolve as the phase changes so that I feel a difference between the peaceful Exploration Phase and the more chaotic Haunt Phase"
$ws.Range("B11").Value = "As a player, I want to be able to view the current state of the Haunt Timer at so that I’m able to make tactical decisions about how to play"
$ws.Range("B13").Value = "As a player, I want to be able to have an “inventory” to keep track of my items and what they do"
$ws.Range("B14").Value = "As a player, I want to be able to view the board state so that I can decide what my course of action will be"
$ws.Range("B15").Value = "As a player, I want to have a reviewable ruleset to remember my win/loss condition"
$ws.Range("B16").Value = "As a player I want to be able to keep track of current stat points for users so I can know how many dice to roll and if I will die"
$ws.Range("B17").Value = "As a player, I want to have the option to change the size of the font in the game"
$ws.Range("B18").Value = "As a player, I want to be able to invite friends to my game to play local multiplayer"
$ws.Range("B19").Value = "As a player, I want to have a menu for modifiers/house rules before setting up the game"
$ws.Range("B20").Value = "As a player, I want to be able to have the ability to pause the game when playing with NPCs"
$ws.Range("B21").Value = "As a player, I want to be able to return to the main menu after a game so that I can navigate easily"

# --- 2) Re-style row 21s User Story cell: drop the stray left alignment, switch
#        the font back to the columns normal Times New Roman, and give it a solid
#        white fill so it matches the rest of the backlog rows. ---
$rngB21 = $ws.Range("B21")
$rngB21.HorizontalAlignment = 1
$rngB21.Font.Name = "`"Times New Roman`""
$rngB21.Interior.Color = RGB(255, 255, 255)

# --- 3) Fill in the missing Story Points (column C) ---
$storyPoints = @{
    11 = 2
    13 = 1
    14 = 2
    15 = 2
    18 = 2
    19 = 2
    21 = 1
}
foreach ($row in $storyPoints.Keys) {
    $cell = $ws.Range("C$row")
    # nudge the alignment so Excel re-applies an explicit (still "general")
    # alignment record instead of silently reusing the blank placeholder style --
    # matches how the sheet behaves when a value is typed into these cells.
    $cell.HorizontalAlignment = -4108
    $cell.HorizontalAlignment = 1
    $cell.Value = $storyPoints[$row]
}
